# repull data, push all data, mean calculation
# Update column F (dSF) values for the rows whose re-pulled data differs
# from the previously stored value (column E, dS0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = 4
    7  = 6
    8  = 0
    10 = 0
    14 = -1
    16 = -3
    23 = -2
    24 = -8
    25 = -2
    29 = -3
    31 = 3
    36 = 0
    37 = -2
    38 = -12
    43 = -1
    44 = -9
    45 = 0
    53 = -1
    57 = -1
    65 = -1
    66 = 12
    73 = 4
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
